$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank row that sits above the data (row 2), which
# shifts the "test terms" data row (row 3 -> row 2) and the blank
# styled rows below it up by one, and drops the now-superfluous last
# blank row (old row 6).
$ws.Rows.Item(2).Delete()

# Update the active selection to match the post-edit cursor position.
$ws.Range("F6").Select()
